{"js": "// Docassemble template syntax update:\n//   X.name.full(middle=\u2019full\u2019)         ->  X.name_full()\n//   {{comma_and_list(property_successors)}}  ->  {{property_successors.full_names())}}\n//   {{comma_and_list(health_successors)}}    ->  {{health_successors.full_names()}}\n//\n// NOTE: occurrences written as \"middle= \u2018full\u2019\" / \"middle = \u2018full\u2019\" (a SPACE\n// before the closing-quote value) are intentionally left untouched \u2014 the\n// source diff does not rewrite those two spots.\n\nasync function replaceAll(context, findText, replaceText) {\n  const results = context.document.body.search(findText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Simple \"X.name.full(middle=\u2019full\u2019)\" -> \"X.name_full()\" swaps.\n// Doing \"property_agent...}\" also rewrites \"new_property_agent...}\" because\n// the shorter needle is a substring of the longer one (same for health).\nconst swaps = [\n  [\"user.name.full(middle=\\u2019full\\u2019)\", \"user.name_full()\"],\n  [\"person.name.full(middle=\\u2019full\\u2019)\", \"person.name_full()\"],\n  [\n    \"property_agent.name.full(middle=\\u2019full\\u2019)\",\n    \"property_agent.name_full()\",\n  ],\n  [\n    \"health_agent.name.full(middle=\\u2019full\\u2019)\",\n    \"health_agent.name_full()\",\n  ],\n];\n\nfor (const [find, replace] of swaps) {\n  await replaceAll(context, find, replace);\n}\n\n// comma_and_list(...) -> .full_names() rewrites (text of the source diff,\n// including the extra trailing \")\" that the original commit leaves on the\n// property_successors line but not on the health_successors line).\nawait replaceAll(\n  context,\n  \"{{comma_and_list(property_successors)}}\",\n  \"{{property_successors.full_names())}}\"\n);\nawait replaceAll(\n  context,\n  \"comma_and_list(health_successors)\",\n  \"health_successors.full_names()\"\n);\n", "ps1": "# Docassemble template syntax update:\n#   X.name.full(middle=\u2019full\u2019)               -> X.name_full()\n#   {{comma_and_list(property_successors)}}  -> {{property_successors.full_names())}}\n#   {{comma_and_list(health_successors)}}    -> {{health_successors.full_names()}}\n#\n# NOTE: occurrences written as \"middle= \u2018full\u2019\" / \"middle = \u2018full\u2019\" (a SPACE\n# before the closing-quote value) are intentionally left untouched \u2014 the\n# source diff does not rewrite those two spots.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.MatchWildcards = $false\n    $find.MatchCase = $true\n    $find.Replacement.Text = $replaceText\n    # FindWhat, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap (wdFindContinue=1), Format,\n    # ReplaceWith, Replace (wdReplaceAll=2)\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# Simple \"X.name.full(middle=\u2019full\u2019)\" -> \"X.name_full()\" swaps.\n# Doing \"property_agent...)\" also rewrites \"new_property_agent...)\" because\n# the shorter needle is a substring of the longer one (same for health).\nReplace-AllText \"user.name.full(middle=\u2019full\u2019)\" \"user.name_full()\"\nReplace-AllText \"person.name.full(middle=\u2019full\u2019)\" \"person.name_full()\"\nReplace-AllText \"property_agent.name.full(middle=\u2019full\u2019)\" \"property_agent.name_full()\"\nReplace-AllText \"health_agent.name.full(middle=\u2019full\u2019)\" \"health_agent.name_full()\"\n\n# comma_and_list(...) -> .full_names() rewrites (text of the source diff,\n# including the extra trailing \")\" that the original commit leaves on the\n# property_successors line but not on the health_successors line).\nReplace-AllText \"{{comma_and_list(property_successors)}}\" \"{{property_successors.full_names())}}\"\nReplace-AllText \"comma_and_list(health_successors)\" \"health_successors.full_names()\"\n"}
